# Edit slide 5 of the unikernel figure:
#  1. Shape "Rectangle 13" ("LibOS + PL Runtime") -> "PL Runtime"
#     (drop the leading "LibOS" run, tighten " + PL " to "PL ").
#  2. Shape "Rectangle 22" ("Application + Database library"):
#       - grow the box (it now needs to fit an extra line of text)
#       - text becomes "Application + Database library + LibOS"
#         split across three runs so "LibOS" keeps its own run
#         (matching the proofing/err-flagged run used elsewhere
#         in the deck for the same word).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- Shape "Rectangle 13": "LibOS + PL Runtime" -> "PL Runtime" ---
$libOsPl = $s.Shapes.Item(5)
$tr1 = $libOsPl.TextFrame.TextRange

# Delete the "LibOS" run entirely (chars 1-5).
$tr1.Characters(1, 5).Text = ""

# What's left is " + PL " (6 chars) followed by "Runtime"; tighten it to "PL ".
$tr1.Characters(1, 6).Text = "PL "

# --- Shape "Rectangle 22": resize + "Application + Database library" -> "...library + LibOS" ---
$appDb = $s.Shapes.Item(12)

# Resize/reposition the box to make room for the extra text.
$appDb.Top = 215.20009
$appDb.Height = 58.4

$tr2 = $appDb.TextFrame.TextRange

# "library" (chars 24-30) becomes "library + ", then append a separate "LibOS" run.
$tr2.Characters(24, 7).Text = "library + "
$tr2.InsertAfter("LibOS") | Out-Null
